$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dst = $ws.Range("Q1")
$dst.Value = "Docente Correo"
Write-Host "set value only"
